$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '43.037.40'

$ws.Range("E2").Value = '  +0.04%  '

$ws.Range("D3").Value = '2.300.50'

$ws.Range("E3").Value = '  -0.02%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '300.18'
$ws.Range("D5").Style = "Normal"

$ws.Range("E5").Value = '  -0.13%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '97.90'
$ws.Range("D6").Style = "Normal"

$ws.Range("E6").Value = '  -1.49%  '

$ws.Range("E8").Value = '  -0.03%  '

$ws.Range("E9").Value = '  +1.26%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '36.14'
$ws.Range("D10").Style = "Normal"

$ws.Range("E10").Value = '  -0.23%  '

$ws.Range("E11").Value = '  +0.32%  '

$ws.Range("E12").Value = '  +0.75%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '17.73'
$ws.Range("D13").Style = "Normal"

$ws.Range("E13").Value = '  -2.42%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.89'
$ws.Range("D14").Style = "Normal"

$ws.Range("E14").Value = '  -0.48%  '

$ws.Range("D15").Value = '2.659.62'

$ws.Range("E15").Value = '  +0.04%  '

$ws.Range("D16").Value = '2.258.05'

$ws.Range("E16").Value = '  -2.61%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.788'
$ws.Range("D17").Style = "Normal"

$ws.Range("E17").Value = '  -1.41%  '

$ws.Range("D18").Value = '42.917.13'

$ws.Range("E18").Value = '  +0.01%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '12.83'
$ws.Range("D19").Style = "Normal"

$ws.Range("E19").Value = '  +1.99%  '

$ws.Range("E20").Value = '  +0.97%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.13'
$ws.Range("D21").Style = "Normal"

$ws.Range("E21").Value = '  +0.17%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '68.26'
$ws.Range("D22").Style = "Normal"

$ws.Range("E22").Value = '  +0.60%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '237.81'
$ws.Range("D23").Style = "Normal"

$ws.Range("E23").Value = '  +1.04%  '

$ws.Range("E24").Value = '  -0.85%  '

$ws.Range("E25").Value = '  -0.51%  '

$ws.Range("E26").Value = '  -0.74%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '4.02'
$ws.Range("D27").Style = "Normal"

$ws.Range("E27").Value = '  -0.27%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '25.00'
$ws.Range("D28").Style = "Normal"

$ws.Range("E28").Value = '  +0.29%  '

$ws.Range("E29").Value = '  -13.10%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '9.15'
$ws.Range("D30").Style = "Normal"

$ws.Range("E30").Value = '  +0.19%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '163.31'
$ws.Range("D31").Style = "Normal"

$ws.Range("E31").Value = '  -2.45%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '33.09'
$ws.Range("D32").Style = "Normal"

$ws.Range("E32").Value = '  -4.02%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.999'
$ws.Range("D33").Style = "Normal"

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.11'
$ws.Range("D34").Style = "Normal"

$ws.Range("E34").Value = '  +1.82%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '18.13'
$ws.Range("D35").Style = "Normal"

$ws.Range("E35").Value = '  +2.90%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '4.78'
$ws.Range("D36").Style = "Normal"

$ws.Range("E36").Value = '  +4.14%  '

$ws.Range("E37").Value = '  +0.31%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.0697'
$ws.Range("D38").Style = "Normal"

$ws.Range("E38").Value = '  +1.13%  '

$ws.Range("E39").Value = '  +1.12%  '

$ws.Range("E40").Value = '  -0.49%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.78'
$ws.Range("D41").Style = "Normal"

$ws.Range("E41").Value = '  -0.85%  '

$ws.Range("E42").Value = '  +1.31%  '

$ws.Range("D43").Value = '2.018.25'

$ws.Range("E43").Value = '  +1.88%  '

$ws.Range("B44").Value = 'VeChain'

$ws.Range("C44").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0286'
$ws.Range("D44").Style = "Normal"

$ws.Range("E44").Value = '  -1.36%  '

$ws.Range("B45").Value = 'ApeXProtocol'

$ws.Range("C45").Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.25'
$ws.Range("D45").Style = "Normal"

$ws.Range("E45").Value = '  -1.27%  '

$ws.Range("E46").Value = '  +1.67%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '17.54'
$ws.Range("D47").Style = "Normal"

$ws.Range("E47").Value = '  +0.71%  '

$ws.Range("E48").Value = '  -2.37%  '

$ws.Range("E49").Value = '  -1.96%  '

$ws.Range("D50").Value = '2.526.46'

$ws.Range("E50").Value = '  +0.04%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.53'
$ws.Range("D51").Style = "Normal"

$ws.Range("E51").Value = '  -1.54%  '
